$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new Job Posting row with Job_Id = JD_004
$ws.Range("A5").Value = "JD_004"
$ws.Range("B5").Value = "Senior RPA Developer"
$ws.Range("C5").Value = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment`n"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 2
$ws.Rows.Item(5).EntireRow.AutoFit()
